$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 829: copy the formatting from the previous header block (row 793) ---
$ws.Range("A793:H793").Copy()
$ws.Range("A829:H829").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("A829").Value = "States/UT"
$ws.Range("B829").Value = "Active Cases"
$ws.Range("C829").Value = "Active Cases Since Yesterday"
$ws.Range("D829").Value = "Recovered Cases"
$ws.Range("E829").Value = "Recovered Cases Since Yesterday"
$ws.Range("F829").Value = "Deceased Cases"
$ws.Range("G829").Value = "Deceased Cases Since Yesterday"
$ws.Range("H829").Value = "Date"

# --- Data rows 830-864: one row per State/UT for 10-10-2020 ---
$states = @(
  @("Andaman and Nicobar Islands", 197, 7, 3724, 17, 55, 0, "10-10-2020"),
  @("Andhra Pradesh", 47665, -996, 691040, 6110, 6159, 31, "10-10-2020"),
  @("Arunachal Pradesh", 2860, 82, 8877, 198, 23, 1, "10-10-2020"),
  @("Assam", 29710, -1057, 161904, 2068, 802, 8, "10-10-2020"),
  @("Bihar", 11274, -173, 182121, 1425, 934, 5, "10-10-2020"),
  @("Chandigarh", 1292, -100, 11505, 161, 188, 2, "10-10-2020"),
  @("Chhattisgarh", 27439, 12, 108935, 2908, 1196, 38, "10-10-2020"),
  @("Dadra and Nagar Haveli and Daman and Diu", 101, -4, 3037, 12, 2, 0, "10-10-2020"),
  @("Delhi", 21955, -277, 276046, 3098, 5692, 39, "10-10-2020"),
  @("Goa", 4783, 67, 32317, 415, 491, 7, "10-10-2020"),
  @("Gujarat", 16181, -284, 129304, 1518, 3547, 9, "10-10-2020"),
  @("Haryana", 10830, -37, 127540, 1273, 1562, 14, "10-10-2020"),
  @("Himachal Pradesh", 2856, -87, 13876, 279, 245, 7, "10-10-2020"),
  @("Jammu and Kashmir", 11144, -338, 69979, 959, 1306, 15, "10-10-2020"),
  @("Jharkhand", 8819, -453, 81654, 1215, 781, 6, "10-10-2020"),
  @("Karnataka", 118870, 1708, 561610, 9091, 9789, 114, "10-10-2020"),
  @("Kerala", 91841, 1177, 175304, 8048, 955, 25, "10-10-2020"),
  @("Ladakh", 1038, -261, 3886, 346, 63, 0, "10-10-2020"),
  @("Madhya Pradesh", 16168, -620, 124887, 2200, 2574, 27, "10-10-2020"),
  @("Maharashtra", 236947, -5491, 1229339, 17323, 39732, 302, "10-10-2020"),
  @("Manipur", 2858, -19, 9866, 147, 86, 3, "10-10-2020"),
  @("Meghalaya", 2424, 55, 4903, 71, 61, 1, "10-10-2020"),
  @("Mizoram", 188, -32, 1974, 37, 0, 0, "10-10-2020"),
  @("Nagaland", 1215, 60, 5656, 92, 17, 0, "10-10-2020"),
  @("Odisha", 25460, -724, 220388, 3404, 991, 17, "10-10-2020"),
  @("Puducherry", 4803, 76, 25543, 287, 558, 2, "10-10-2020"),
  @("Punjab", 10153, -622, 108533, 1333, 3773, 32, "10-10-2020"),
  @("Rajasthan", 21398, 16, 131766, 2148, 1621, 16, "10-10-2020"),
  @("Sikkim", 504, -41, 2721, 71, 53, 2, "10-10-2020"),
  @("Tamil Nadu", 44197, -240, 591811, 5357, 10120, 68, "10-10-2020"),
  @("Telengana", 26104, -270, 183025, 2072, 1217, 9, "10-10-2020"),
  @("Tripura", 4039, -158, 23801, 327, 313, 2, "10-10-2020"),
  @("Uttarakhand", 7289, -560, 46058, 1250, 716, 14, "10-10-2020"),
  @("Uttar Pradesh", 41287, -1265, 383086, 4424, 6293, 48, "10-10-2020"),
  @("West Bengal", 29296, 442, 252806, 3069, 5501, 62, "10-10-2020")
)

# Column H holds day-month-year text like "10-10-2020" which Excel would
# otherwise auto-convert to a date serial; format the destination cells as
# Text first so the literal string is preserved, matching the rest of the
# sheet (which stores these as plain text).
$ws.Range("H830:H864").NumberFormat = "@"

$row = 830
foreach ($s in $states) {
  $ws.Cells.Item($row, 1).Value = $s[0]
  $ws.Cells.Item($row, 2).Value = $s[1]
  $ws.Cells.Item($row, 3).Value = $s[2]
  $ws.Cells.Item($row, 4).Value = $s[3]
  $ws.Cells.Item($row, 5).Value = $s[4]
  $ws.Cells.Item($row, 6).Value = $s[5]
  $ws.Cells.Item($row, 7).Value = $s[6]
  $ws.Cells.Item($row, 8).Value = $s[7]
  $row++
}

# The rest of the sheet's data rows carry no explicit style (style index 0),
# so drop the temporary Text number format now that the literal strings are
# safely stored, to match the unstyled look of every other data row.
$ws.Range("H830:H864").Style = "Normal"
